$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 194; this pushes the existing rows 194-220
# down to 196-222 and extends the sheet dimension to A1:T222.
$ws.Range("A194:A195").EntireRow.Insert()

# Populate the two newly inserted rows (194 and 195) with the new weekly
# price observations (same market/product context as the surrounding rows).
$ws.Cells.Item(194, 1).Value2 = 4
$ws.Cells.Item(194, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(194, 3).Value2 = "Los Lagos"
$ws.Cells.Item(194, 4).Value2 = 44449
$ws.Cells.Item(194, 5).Value2 = 10
$ws.Cells.Item(194, 6).Value2 = "Fruta"
$ws.Cells.Item(194, 7).Value2 = 100102
$ws.Cells.Item(194, 8).Value2 = "Cítricos"
$ws.Cells.Item(194, 9).Value2 = 100102005
$ws.Cells.Item(194, 10).Value2 = "Naranja"
$ws.Cells.Item(194, 11).Value2 = "Navel Late"
$ws.Cells.Item(194, 12).Value2 = "Primera"
$ws.Cells.Item(194, 13).Value2 = 400
$ws.Cells.Item(194, 14).Value2 = 13000
$ws.Cells.Item(194, 15).Value2 = 13000
$ws.Cells.Item(194, 16).Value2 = 13000
$ws.Cells.Item(194, 17).Value2 = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(194, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(194, 19).Value2 = 867
$ws.Cells.Item(194, 20).Value2 = 15

$ws.Cells.Item(195, 1).Value2 = 4
$ws.Cells.Item(195, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(195, 3).Value2 = "Los Lagos"
$ws.Cells.Item(195, 4).Value2 = 44449
$ws.Cells.Item(195, 5).Value2 = 10
$ws.Cells.Item(195, 6).Value2 = "Fruta"
$ws.Cells.Item(195, 7).Value2 = 100102
$ws.Cells.Item(195, 8).Value2 = "Cítricos"
$ws.Cells.Item(195, 9).Value2 = 100102005
$ws.Cells.Item(195, 10).Value2 = "Naranja"
$ws.Cells.Item(195, 11).Value2 = "Navel Late"
$ws.Cells.Item(195, 12).Value2 = "Segunda"
$ws.Cells.Item(195, 13).Value2 = 200
$ws.Cells.Item(195, 14).Value2 = 10000
$ws.Cells.Item(195, 15).Value2 = 10000
$ws.Cells.Item(195, 16).Value2 = 10000
$ws.Cells.Item(195, 17).Value2 = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(195, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(195, 19).Value2 = 667
$ws.Cells.Item(195, 20).Value2 = 15
